{"js": "// Replace the date line and each \"NNN\u00d7N=NNNN\" multiplication-table entry\n// with its new value, matching the unified diff 1:1 (each original text\n// string is unique in the document, so a plain search+replace is safe).\nconst replacements = [\n  [\"2024-07-06 Saturday\", \"2024-07-07 Sunday\"],\n  [\"835\u00d79=7515\", \"187\u00d77=1309\"],\n  [\"822\u00d75=4110\", \"465\u00d73=1395\"],\n  [\"340\u00d76=2040\", \"845\u00d74=3380\"],\n  [\"216\u00d76=1296\", \"711\u00d76=4266\"],\n  [\"369\u00d76=2214\", \"256\u00d74=1024\"],\n  [\"146\u00d78=1168\", \"251\u00d78=2008\"],\n  [\"749\u00d79=6741\", \"792\u00d74=3168\"],\n  [\"971\u00d72=1942\", \"982\u00d72=1964\"],\n  [\"556\u00d76=3336\", \"703\u00d77=4921\"],\n  [\"474\u00d76=2844\", \"577\u00d78=4616\"],\n  [\"751\u00d73=2253\", \"790\u00d75=3950\"],\n  [\"411\u00d73=1233\", \"277\u00d73=831\"],\n  [\"386\u00d73=1158\", \"297\u00d79=2673\"],\n  [\"248\u00d79=2232\", \"922\u00d79=8298\"],\n  [\"357\u00d76=2142\", \"185\u00d74=740\"],\n  [\"601\u00d76=3606\", \"345\u00d79=3105\"],\n  [\"823\u00d73=2469\", \"865\u00d75=4325\"],\n  [\"975\u00d79=8775\", \"689\u00d72=1378\"],\n  [\"764\u00d73=2292\", \"412\u00d77=2884\"],\n  [\"781\u00d73=2343\", \"862\u00d77=6034\"],\n  [\"439\u00d72=878\", \"488\u00d75=2440\"],\n  [\"678\u00d74=2712\", \"586\u00d72=1172\"],\n  [\"120\u00d78=960\", \"744\u00d76=4464\"],\n  [\"602\u00d73=1806\", \"662\u00d76=3972\"],\n  [\"713\u00d76=4278\", \"280\u00d77=1960\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"NNN\u00d7N=NNNN\" multiplication-table entry\n# with its new value, matching the unified diff 1:1 (each original text\n# string is unique in the document, so a plain find/replace is safe).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-07-06 Saturday\", \"2024-07-07 Sunday\"),\n    @(\"835\u00d79=7515\", \"187\u00d77=1309\"),\n    @(\"822\u00d75=4110\", \"465\u00d73=1395\"),\n    @(\"340\u00d76=2040\", \"845\u00d74=3380\"),\n    @(\"216\u00d76=1296\", \"711\u00d76=4266\"),\n    @(\"369\u00d76=2214\", \"256\u00d74=1024\"),\n    @(\"146\u00d78=1168\", \"251\u00d78=2008\"),\n    @(\"749\u00d79=6741\", \"792\u00d74=3168\"),\n    @(\"971\u00d72=1942\", \"982\u00d72=1964\"),\n    @(\"556\u00d76=3336\", \"703\u00d77=4921\"),\n    @(\"474\u00d76=2844\", \"577\u00d78=4616\"),\n    @(\"751\u00d73=2253\", \"790\u00d75=3950\"),\n    @(\"411\u00d73=1233\", \"277\u00d73=831\"),\n    @(\"386\u00d73=1158\", \"297\u00d79=2673\"),\n    @(\"248\u00d79=2232\", \"922\u00d79=8298\"),\n    @(\"357\u00d76=2142\", \"185\u00d74=740\"),\n    @(\"601\u00d76=3606\", \"345\u00d79=3105\"),\n    @(\"823\u00d73=2469\", \"865\u00d75=4325\"),\n    @(\"975\u00d79=8775\", \"689\u00d72=1378\"),\n    @(\"764\u00d73=2292\", \"412\u00d77=2884\"),\n    @(\"781\u00d73=2343\", \"862\u00d77=6034\"),\n    @(\"439\u00d72=878\", \"488\u00d75=2440\"),\n    @(\"678\u00d74=2712\", \"586\u00d72=1172\"),\n    @(\"120\u00d78=960\", \"744\u00d76=4464\"),\n    @(\"602\u00d73=1806\", \"662\u00d76=3972\"),\n    @(\"713\u00d76=4278\", \"280\u00d77=1960\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
